# Regen save_data to use K instead of Strike#: update column G (K) values
# for the affected rows on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 2
    9  = 2
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 2
    17 = 2
    19 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
